$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'69.491.93"
$ws.Range("D2").Style = $style
$style = $ws.Range("E2").Style
$ws.Range("E2").Value = "'  +2.14%  "
$ws.Range("E2").Style = $style
$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.432.31"
$ws.Range("D3").Style = $style
$style = $ws.Range("E3").Style
$ws.Range("E3").Value = "'  +0.90%  "
$ws.Range("E3").Style = $style
$style = $ws.Range("E4").Style
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = $style
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'564.47"
$ws.Range("D5").Style = $style
$style = $ws.Range("E5").Style
$ws.Range("E5").Value = "'  +2.15%  "
$ws.Range("E5").Style = $style
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'167.19"
$ws.Range("D6").Style = $style
$style = $ws.Range("E6").Style
$ws.Range("E6").Value = "'  +5.40%  "
$ws.Range("E6").Style = $style
$style = $ws.Range("E7").Style
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = $style
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.514"
$ws.Range("D8").Style = $style
$style = $ws.Range("E8").Style
$ws.Range("E8").Value = "'  +1.92%  "
$ws.Range("E8").Style = $style
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.170"
$ws.Range("D9").Style = $style
$style = $ws.Range("E9").Style
$ws.Range("E9").Value = "'  +7.45%  "
$ws.Range("E9").Style = $style
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'2.428.98"
$ws.Range("D10").Style = $style
$style = $ws.Range("E10").Style
$ws.Range("E10").Value = "'  +1.37%  "
$ws.Range("E10").Style = $style
$style = $ws.Range("E11").Style
$ws.Range("E11").Value = "'  -1.77%  "
$ws.Range("E11").Style = $style
$style = $ws.Range("E12").Style
$ws.Range("E12").Value = "'  +2.07%  "
$ws.Range("E12").Style = $style
$style = $ws.Range("E13").Style
$ws.Range("E13").Value = "'  -1.08%  "
$ws.Range("E13").Style = $style
$style = $ws.Range("E14").Style
$ws.Range("E14").Value = "'  +5.74%  "
$ws.Range("E14").Style = $style
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'69.291.31"
$ws.Range("D15").Style = $style
$style = $ws.Range("E15").Style
$ws.Range("E15").Value = "'  +2.00%  "
$ws.Range("E15").Style = $style
$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'2.877.19"
$ws.Range("D16").Style = $style
$style = $ws.Range("E16").Style
$ws.Range("E16").Value = "'  -0.98%  "
$ws.Range("E16").Style = $style
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'23.97"
$ws.Range("D17").Style = $style
$style = $ws.Range("E17").Style
$ws.Range("E17").Value = "'  +4.67%  "
$ws.Range("E17").Style = $style
$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'2.424.83"
$ws.Range("D18").Style = $style
$style = $ws.Range("E18").Style
$ws.Range("E18").Value = "'  +0.79%  "
$ws.Range("E18").Style = $style
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'10.83"
$ws.Range("D19").Style = $style
$style = $ws.Range("E19").Style
$ws.Range("E19").Value = "'  +5.06%  "
$ws.Range("E19").Style = $style
$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'343.44"
$ws.Range("D20").Style = $style
$style = $ws.Range("E20").Style
$ws.Range("E20").Value = "'  +4.36%  "
$ws.Range("E20").Style = $style
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'7.15"
$ws.Range("D21").Style = $style
$style = $ws.Range("E21").Style
$ws.Range("E21").Value = "'  +5.81%  "
$ws.Range("E21").Style = $style
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'3.87"
$ws.Range("D22").Style = $style
$style = $ws.Range("E22").Style
$ws.Range("E22").Value = "'  +2.92%  "
$ws.Range("E22").Style = $style
$style = $ws.Range("E23").Style
$ws.Range("E23").Value = "'  +6.29%  "
$ws.Range("E23").Style = $style
$style = $ws.Range("E24").Style
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("E24").Style = $style
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'66.08"
$ws.Range("D25").Style = $style
$style = $ws.Range("E25").Style
$ws.Range("E25").Value = "'  +0.62%  "
$ws.Range("E25").Style = $style
$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'3.82"
$ws.Range("D26").Style = $style
$style = $ws.Range("E26").Style
$ws.Range("E26").Value = "'  +5.66%  "
$ws.Range("E26").Style = $style
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'1.03"
$ws.Range("D27").Style = $style
$style = $ws.Range("E27").Style
$ws.Range("E27").Value = "'  +3.24%  "
$ws.Range("E27").Style = $style
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'8.49"
$ws.Range("D28").Style = $style
$style = $ws.Range("E28").Style
$ws.Range("E28").Value = "'  +5.98%  "
$ws.Range("E28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'2.554.20"
$ws.Range("D29").Style = $style
$style = $ws.Range("E29").Style
$ws.Range("E29").Value = "'  +1.18%  "
$ws.Range("E29").Style = $style
$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.0₃0853"
$ws.Range("D30").Style = $style
$style = $ws.Range("E30").Style
$ws.Range("E30").Value = "'  +7.17%  "
$ws.Range("E30").Style = $style
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'7.40"
$ws.Range("D31").Style = $style
$style = $ws.Range("E31").Style
$ws.Range("E31").Value = "'  +5.32%  "
$ws.Range("E31").Style = $style
$style = $ws.Range("E32").Style
$ws.Range("E32").Value = "'  +10.72%  "
$ws.Range("E32").Style = $style
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'455.70"
$ws.Range("D33").Style = $style
$style = $ws.Range("E33").Style
$ws.Range("E33").Value = "'  +8.27%  "
$ws.Range("E33").Style = $style
$style = $ws.Range("E34").Style
$ws.Range("E34").Value = "'  +0.12%  "
$ws.Range("E34").Style = $style
$style = $ws.Range("E35").Style
$ws.Range("E35").Value = "'  +2.63%  "
$ws.Range("E35").Style = $style
$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'158.93"
$ws.Range("D36").Style = $style
$style = $ws.Range("E36").Style
$ws.Range("E36").Value = "'  +0.36%  "
$ws.Range("E36").Style = $style
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'19.10"
$ws.Range("D37").Style = $style
$style = $ws.Range("E37").Style
$ws.Range("E37").Value = "'  +0.57%  "
$ws.Range("E37").Style = $style
$style = $ws.Range("E38").Style
$ws.Range("E38").Value = "'  +6.26%  "
$ws.Range("E38").Style = $style
$style = $ws.Range("E39").Style
$ws.Range("E39").Value = "'  +0.00%  "
$ws.Range("E39").Style = $style
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'18.25"
$ws.Range("D40").Style = $style
$style = $ws.Range("E40").Style
$ws.Range("E40").Value = "'  +3.58%  "
$ws.Range("E40").Style = $style
$style = $ws.Range("E41").Style
$ws.Range("E41").Value = "'  +3.93%  "
$ws.Range("E41").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'4.42"
$ws.Range("D42").Style = $style
$style = $ws.Range("E42").Style
$ws.Range("E42").Value = "'  +3.94%  "
$ws.Range("E42").Style = $style
$style = $ws.Range("E43").Style
$ws.Range("E43").Value = "'  +4.68%  "
$ws.Range("E43").Style = $style
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'37.84"
$ws.Range("D44").Style = $style
$style = $ws.Range("E44").Style
$ws.Range("E44").Value = "'  +1.96%  "
$ws.Range("E44").Style = $style
$style = $ws.Range("E45").Style
$ws.Range("E45").Value = "'  +2.85%  "
$ws.Range("E45").Style = $style
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.10"
$ws.Range("D46").Style = $style
$style = $ws.Range("E46").Style
$ws.Range("E46").Value = "'  +7.80%  "
$ws.Range("E46").Style = $style
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'136.26"
$ws.Range("D47").Style = $style
$style = $ws.Range("E47").Style
$ws.Range("E47").Value = "'  +6.09%  "
$ws.Range("E47").Style = $style
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'3.40"
$ws.Range("D48").Style = $style
$style = $ws.Range("E48").Style
$ws.Range("E48").Value = "'  +3.63%  "
$ws.Range("E48").Style = $style
$style = $ws.Range("E49").Style
$ws.Range("E49").Value = "'  +2.78%  "
$ws.Range("E49").Style = $style
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.490"
$ws.Range("D50").Style = $style
$style = $ws.Range("E50").Style
$ws.Range("E50").Value = "'  +3.68%  "
$ws.Range("E50").Style = $style
$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'0.0935"
$ws.Range("D51").Style = $style
$style = $ws.Range("E51").Style
$ws.Range("E51").Value = "'  +2.57%  "
$ws.Range("E51").Style = $style
